$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell while preserving its
# original (text) storage type and original cell style. Assigning a plain
# numeric-looking string via .Value would otherwise get auto-converted to a
# real number by Excel; temporarily forcing a text number format avoids that,
# and we restore the original style afterwards so the cell's formatting is
# left exactly as it was.
function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

# Enterprises density (per 1000 people) - Source Type: Statistical Institution
Set-TextValue "B11" "34.32"
Set-TextValue "C11" "7.02"
Set-TextValue "D11" "41.34"

# Employment (% of total) - Source Type: Statistical Institution
Set-TextValue "B12" "14.44"
Set-TextValue "C12" "26.57"
Set-TextValue "D12" "41.01"

# Enterprises density (per 1000 people) - Source Type: SME Associations
Set-TextValue "B33" "24.19"
Set-TextValue "C33" "2.75"
Set-TextValue "D33" "26.94"

# Employment (% of total) - Source Type: SME Associations
Set-TextValue "B34" "17.15"
Set-TextValue "C34" "35.95"

# Enterprises (% of total)
Set-TextValue "B36" "89.49"
Set-TextValue "C36" "10.17"
Set-TextValue "D36" "99.65"

# Value added to the economy (% of total)
Set-TextValue "B40" "18.66"
Set-TextValue "C40" "32.39"
Set-TextValue "D40" "51.05"
